$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 35) below the existing data, matching the
# date formatting used by the rest of the Date column.
$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A35").Value = 46006
$ws.Range("B35").Value = 3

# Update the visible selection to match the new last row (mirrors Excel
# auto-selecting/scrolling to the newly entered row).
$ws.Range("A35:B35").Select()
$excel.ActiveWindow.ScrollRow = 19
